# Fruta / hortaliza, semanal
# New weekly data point added at the top of this product's history block.
# This inserts a new row at row 400, pushing the existing rows 400-463
# down to 401-464, and populates the new row 400 with the latest record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 400 (shifts rows 400:463 -> 401:464)
$ws.Rows(400).Insert()

# Populate the newly inserted row 400 with the new weekly observation
$ws.Range("A400").Value = 4
$ws.Range("B400").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C400").Value = "Los Lagos"
$ws.Range("D400").Value = 45180
$ws.Range("E400").Value = 10
$ws.Range("F400").Value = "Fruta"
$ws.Range("G400").Value = 100108
$ws.Range("H400").Value = "Tropicales y subtropicales"
$ws.Range("I400").Value = 100108005
$ws.Range("J400").Value = "Piña"
$ws.Range("K400").Value = "Caramelo"
$ws.Range("L400").Value = "Primera"
$ws.Range("M400").Value = 25
$ws.Range("N400").Value = 25000
$ws.Range("O400").Value = 25000
$ws.Range("P400").Value = 25000
$ws.Range("Q400").Value = "$/caja 14 unidades"
$ws.Range("R400").Value = "Ecuador"
$ws.Range("S400").Value = 1786
$ws.Range("T400").Value = 14
